$d = $word.ActiveDocument

# --- Text replacements (find unique paragraph text, set Range.Text directly
# to avoid smart-quote autocorrection during Find/Replace) ---
function Replace-ParagraphText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $newText
    }
    return $found
}

$null = Replace-ParagraphText "WhatsApp-konversation mellan Heba Alhussien och ""Abu Mohammad"" (professionell penningväxlare/صراف) under perioden december 2023 – september 2024." "WhatsApp-konversation mellan Heba och Abu Mohammad (penningväxlare/صراف), dec 2023–sep 2024."
$null = Replace-ParagraphText "Att Heba bedrev systematisk hawala-verksamhet med professionella penningväxlare" "Heba bedrev systematisk hawala med professionella penningväxlare"
$null = Replace-ParagraphText "Att pengarna som påstås vara ""lån"" i själva verket var del av denna hawala-verksamhet" "Påstådda ""lån"" var del av denna verksamhet"
$null = Replace-ParagraphText "Att identiskt transaktionsmönster (fakturor, Swish, dollarkurser) förekommer med Abu Mohammad" "Identiskt transaktionsmönster (fakturor, Swish, dollarkurser) som med Mohammad"
$null = Replace-ParagraphText "518 meddelanden" "518 meddelanden, 60+ ljud, 80+ foton, 15+ PDF-fakturor"
$null = Replace-ParagraphText "Dollarväxlingar i Damaskus" "Dollarväxlingar, 11,5 milj. SYP, USDT/krypto"
$null = Replace-ParagraphText "Hur dollarkurser och belopp fastställdes" "Dollarkurser, betalningsflöden och internationella transaktioner"
$null = Replace-ParagraphText "I enlighet med rättens anvisningar kommer samtlig arabisk bevisning att förses med auktoriserad svensk översättning. Översättningarna beställs och levereras före fristen 2026-03-20." "I enlighet med rättens anvisningar förses all arabisk bevisning med auktoriserad svensk översättning före fristen 2026-03-20."

# --- Whole-paragraph deletions ---
# Resolve each unique bullet's paragraph index first, then delete
# from highest index to lowest so earlier indices stay valid.
$deleteTexts = @(
    "60+ ljudmeddelanden med transaktionsinstruktioner",
    "80+ fotografier med kvitton och betalningsbevis",
    "15+ PDF-fakturor (Klarna, Nordea Ropo Capital)",
    "Transaktioner i syriska pund (11,5 miljoner SYP)",
    "USDT/kryptovalutatransaktioner",
    "Att Heba använde sitt svenska konto för internationella transaktioner"
)

$indices = @()
foreach ($txt in $deleteTexts) {
    $rng = $d.Content
    $found = $rng.Find.Execute($txt, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $indices += $rng.Paragraphs.Item(1).Index
    }
}

$indices = $indices | Sort-Object -Descending
foreach ($idx in $indices) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

Write-Output "Done"
